# Pass1_게임서버_4주_학습계획.xlsx - Week2 Day9 chat_client analysis update
#
# - Mark Day 7 (chat_server 분석 2) and Day 8 (chat_server 분석 3) rows as complete
#   in the "완료" (Done) column, matching the checked-box style already used by
#   the other completed rows (e.g. Day 4-6).
# - Move the active cell selection to J13.
# - Set the page setup to A4 / portrait (adds a pageSetup part to the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 8 and 9 (Day 7 / Day 8) as done, reusing the same formatting (checked
# box symbol font, style index 2 in the original file) already applied to the
# previous completed rows (e.g. G5:G7). Copy/PasteSpecial(formats) reuses the
# existing cell style instead of registering a new (duplicate) one.
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G8").Value = "☑"
$ws.Range("G9").Value = "☑"

# Update the selected / active cell.
$ws.Range("J13").Select()

# Configure page setup (A4, portrait) which also materializes a pageSetup part.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
